$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original "Hours", "Module Code" and "Module Title" columns (C,D,E)
# plus the original "Room", "Block" and "Group" columns (H,I,J) for rows 2-10,
# before they get overwritten, since the new layout re-shuffles them.
$rowCount = 10
$origHours  = @{}
$origCode   = @{}
$origTitle  = @{}
$origRoom   = @{}
$origBlock  = @{}
$origGroup  = @{}

for ($r = 2; $r -le $rowCount; $r++) {
    $origHours[$r] = $ws.Cells.Item($r, 5).Value2   # E - Hours
    $origCode[$r]  = $ws.Cells.Item($r, 3).Value2   # C - Module Code
    $origTitle[$r] = $ws.Cells.Item($r, 4).Value2   # D - Module Title
    $origRoom[$r]  = $ws.Cells.Item($r, 8).Value2   # H - Room
    $origBlock[$r] = $ws.Cells.Item($r, 9).Value2   # I - Block
    $origGroup[$r] = $ws.Cells.Item($r, 10).Value2  # J - Group
}

# Row 1 becomes a single title cell; drop the rest of the header labels.
$ws.Range("A1").Value = "Herald College Kathmandu"
$ws.Range("B1:L1").Clear()

# Rewrite rows 2-10 with the new column order:
# A Day | B Time | C Hours | D Module Code | E Module Title | F Class Type |
# G Lecturer | H Group | I Block | J Room
for ($r = 2; $r -le $rowCount; $r++) {
    $ws.Cells.Item($r, 3).Value  = $origHours[$r]   # C - Hours (numeric)
    $ws.Cells.Item($r, 4).Value  = $origCode[$r]    # D - Module Code
    $ws.Cells.Item($r, 5).Value  = $origTitle[$r]   # E - Module Title
    $ws.Cells.Item($r, 8).Value  = $origGroup[$r]   # H - Group
    $ws.Cells.Item($r, 9).Value  = $origBlock[$r]   # I - Block
    $ws.Cells.Item($r, 10).Value = $origRoom[$r]    # J - Room
}

# Remove the now-unused Level and Course columns (K, L).
$ws.Range("K1:L10").Clear()

$ws.Range("A1:J10").Select()
